# Shift each data row (B:K staircase block) one column to the right,
# inserting a new first value in column B (per commit: "Added filtering
# options for the Component Analysis" -> adds a new leading Q-period column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:K2").Value = $ws.Range("B2:J2").Value2
$ws.Range("B2").Value = 0.3332190829615296

$ws.Range("C3:K3").Value = $ws.Range("B3:J3").Value2
$ws.Range("B3").Value = 1.379959960477767

$ws.Range("C4:K4").Value = $ws.Range("B4:J4").Value2
$ws.Range("B4").Value = 0.6848858907743085

$ws.Range("C5:K5").Value = $ws.Range("B5:J5").Value2
$ws.Range("B5").Value = 0.2394483440026288

$ws.Range("C6:K6").Value = $ws.Range("B6:J6").Value2
$ws.Range("B6").Value = 0.8264940743873155

$ws.Range("C7:K7").Value = $ws.Range("B7:J7").Value2
$ws.Range("B7").Value = -0.6581502397256208

$ws.Range("C8:K8").Value = $ws.Range("B8:J8").Value2
$ws.Range("B8").Value = -0.0614365421215774

$ws.Range("C9:K9").Value = $ws.Range("B9:J9").Value2
$ws.Range("B9").Value = 0.9788013920790769

$ws.Range("C10:K10").Value = $ws.Range("B10:J10").Value2
$ws.Range("B10").Value = -0.6334597294260885

$ws.Range("C11:K11").Value = $ws.Range("B11:J11").Value2
$ws.Range("B11").Value = -0.1663369546881469

$ws.Range("C12:K12").Value = $ws.Range("B12:J12").Value2
$ws.Range("B12").Value = 0.04723373655514029

$ws.Range("C13:K13").Value = $ws.Range("B13:J13").Value2
$ws.Range("B13").Value = 0.2540520351237201

$ws.Range("C14:K14").Value = $ws.Range("B14:J14").Value2
$ws.Range("B14").Value = -0.2357884485866682

$ws.Range("C15:K15").Value = $ws.Range("B15:J15").Value2
$ws.Range("B15").Value = -0.01056053025932102

$ws.Range("C16:K16").Value = $ws.Range("B16:J16").Value2
$ws.Range("B16").Value = 0.2172088994749047

$ws.Range("C17:K17").Value = $ws.Range("B17:J17").Value2
$ws.Range("B17").Value = 0.3096861692580615

$ws.Range("C18:K18").Value = $ws.Range("B18:J18").Value2
$ws.Range("B18").Value = -0.154304133832004

$ws.Range("C19:K19").Value = $ws.Range("B19:J19").Value2
$ws.Range("B19").Value = 0.6504264212191833

$ws.Range("C20:K20").Value = $ws.Range("B20:J20").Value2
$ws.Range("B20").Value = -0.2161650486182091

$ws.Range("C21:K21").Value = $ws.Range("B21:J21").Value2
$ws.Range("B21").Value = -0.4070291290349564

$ws.Range("C22:K22").Value = $ws.Range("B22:J22").Value2
$ws.Range("B22").Value = 0.4990422171774198

$ws.Range("C23:K23").Value = $ws.Range("B23:J23").Value2
$ws.Range("B23").Value = -0.1588489131555126

$ws.Range("C24:K24").Value = $ws.Range("B24:J24").Value2
$ws.Range("B24").Value = 0.05616382097024405

$ws.Range("C25:K25").Value = $ws.Range("B25:J25").Value2
$ws.Range("B25").Value = 2.249802839611392

$ws.Range("C26:K26").Value = $ws.Range("B26:J26").Value2
$ws.Range("B26").Value = 7.469150330857293

$ws.Range("C27:K27").Value = $ws.Range("B27:J27").Value2
$ws.Range("B27").Value = -18.17126180013747

$ws.Range("C28:K28").Value = $ws.Range("B28:J28").Value2
$ws.Range("B28").Value = 8.13967154697915

$ws.Range("C29:K29").Value = $ws.Range("B29:J29").Value2
$ws.Range("B29").Value = 2.147570471799392

$ws.Range("C30:K30").Value = $ws.Range("B30:J30").Value2
$ws.Range("B30").Value = -3.44002998652333

$ws.Range("C31:K31").Value = $ws.Range("B31:J31").Value2
$ws.Range("B31").Value = -0.06117417841420103

$ws.Range("C32:K32").Value = $ws.Range("B32:J32").Value2
$ws.Range("B32").Value = 2.041455937156254

$ws.Range("C33:K33").Value = $ws.Range("B33:J33").Value2
$ws.Range("B33").Value = -0.5825945370336409

$ws.Range("C34:K34").Value = $ws.Range("B34:J34").Value2
$ws.Range("B34").Value = 0.09567504080935779

$ws.Range("C35:K35").Value = $ws.Range("B35:J35").Value2
$ws.Range("B35").Value = -0.2604190369987228

$ws.Range("C36:K36").Value = $ws.Range("B36:J36").Value2
$ws.Range("B36").Value = 0.8354549961584912

$ws.Range("C37:K37").Value = $ws.Range("B37:J37").Value2
$ws.Range("B37").Value = -0.1000793599026215

$ws.Range("C38:J38").Value = $ws.Range("B38:I38").Value2
$ws.Range("B38").Value = -0.3537865060796963

$ws.Range("C39:I39").Value = $ws.Range("B39:H39").Value2
$ws.Range("B39").Value = 0.1481773904324453

$ws.Range("C40:H40").Value = $ws.Range("B40:G40").Value2
$ws.Range("B40").Value = 0.157445989004155

$ws.Range("C41:G41").Value = $ws.Range("B41:F41").Value2
$ws.Range("B41").Value = -0.5006594565260708

$ws.Range("C42:F42").Value = $ws.Range("B42:E42").Value2
$ws.Range("B42").Value = 0.2803578805354692

$ws.Range("C43:E43").Value = $ws.Range("B43:D43").Value2
$ws.Range("B43").Value = -0.1719748578450117

$ws.Range("C44:D44").Value = $ws.Range("B44:C44").Value2
$ws.Range("B44").Value = 0.3058625397463315

$ws.Range("C45:C45").Value = $ws.Range("B45:B45").Value2
$ws.Range("B45").Value = -0.6123299526872862

$ws.Range("B46").Value = 0.6883713851991116

$ws.Range("B47").Value = -0.2766911554241067

